$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (renamed headers) ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "use_sublocation"
$ws.Range("C1").Value = "sublocation_min"
$ws.Range("D1").Value = "sublocation_max"

# --- Data rows (values unchanged, but fill in the missing stock room 2 cells) ---
$ws.Range("A2").Value = "active parts"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 100

$ws.Range("A3").Value = "passive parts"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 100

$ws.Range("A4").Value = "stock room 1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 10

$ws.Range("A5").Value = "stock room 2"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# --- Formatting: center-align the value columns (B:D) ---
$ws.Range("B1:D5").HorizontalAlignment = -4108

# --- Column widths (best-fit sizes to match the new, longer headers) ---
$ws.Columns.Item(1).ColumnWidth = 13.42
$ws.Columns.Item(2).ColumnWidth = 14.59
$ws.Columns.Item(3).ColumnWidth = 19.59
$ws.Columns.Item(4).ColumnWidth = 19.75
